$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update existing employee record (kavin -> KAVIN) ---
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = "KAVIN"
# Phone numbers must stay text even though they look numeric (leading zeros
# matter) - route the literal through TEXT()+PasteSpecial so the cell keeps
# the default (unstyled) cell format instead of gaining a new number format.
$ws.Range("Z1").Formula = '=TEXT(90876546,"0")'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4163) | Out-Null
$ws.Range("E2").Value = "JKSDF"

# --- Row 3: new employee record ---
$ws.Range("A3").Value = 4
$ws.Range("B3").Value = "PRAVIN"
$ws.Range("C3").Value = 36489
$ws.Range("Z1").Formula = '=TEXT(89765467,"000000000")'
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D3").PasteSpecial(-4163) | Out-Null
$ws.Range("E3").Value = "SDFKJ"
$ws.Range("F3").Value = 25

# Clean up the scratch cell used to build the text-number values.
$ws.Range("Z1").ClearContents()
$excel.CutCopyMode = 0

Write-Output "done"
